# edit.ps1 - applies the Memoria.docx change:
#  1) Adds a new "1.2) Establecer IP estática..." paragraph right after the
#     VNC paragraph (Raspberry section, item 1).
#  2) Splits the run that used to hold "5) ... web." + the _GoBack bookmark
#     into its own clean paragraph, and appends a brand-new "6) Descargar
#     el proyecto en raspberry mediante Git. ..." paragraph (with spell-check
#     proofErr markers) carrying the _GoBack bookmark at its end.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: insert "1.2) Establecer IP estática ..." paragraph after the
# paragraph that contains "al enchufarla." (end of Raspberry item 1).
# ---------------------------------------------------------------------
$anchor1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*al enchufarla.*") {
        $anchor1 = $p
    }
}

$insertPoint1 = $d.Range($anchor1.Range.End, $anchor1.Range.End)
$xml1 = '<w:p><w:r><w:tab/><w:t>1.2) Establecer IP estática para que siempre sea 192.168.1.137</w:t></w:r></w:p><w:p/>'
$insertPoint1.InsertXML($xml1)

# InsertXML above leaves a stray empty paragraph behind it (an artefact of
# how the host splits multi-<w:p> payloads) - find & remove it.
$stray = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*1.2) Establecer IP estática*") {
        $stray = $p.Next()
    }
}
if ($stray -ne $null -and $stray.Range.Text -eq "") {
    $stray.Range.Delete()
}

# ---------------------------------------------------------------------
# Edit 2: turn the trailing "5) ... web." paragraph (which also carries the
# hidden _GoBack bookmark) into a clean paragraph, and append the new
# "6) Descargar el proyecto ..." paragraph with the bookmark moved to its end.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tailRange = $d.Range($lastPara.Range.Start, $d.Content.End)

$xml2 = '<w:p><w:r><w:t>5) Se consigue encender y apagar un led mediante una web.</w:t></w:r></w:p>' + `
        '<w:p>' + `
        '<w:r><w:t xml:space="preserve">6) Descargar el proyecto en </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:t>raspberry</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> mediante </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve">. Trabajo en el portátil y subo archivos con </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:t>fileZilla</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t>.</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
        '</w:p>'

$tailRange.InsertXML($xml2)
